$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.525.12"
$ws.Range("E2").Value = "  +0.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.480.96"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.22"
$ws.Range("E5").Value = "  +1.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.22"
$ws.Range("E6").Value = "  +0.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.542"
$ws.Range("E7").Value = "  -1.41%  "

$ws.Range("E8").Value = "  -0.12%  "

$ws.Range("E9").Value = "  +2.83%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.73"
$ws.Range("E10").Value = "  -1.25%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("E12").Value = "  +2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.866.00"
$ws.Range("E13").Value = "  +0.81%  "

$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.99"
$ws.Range("E15").Value = "  +9.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.498.14"
$ws.Range("E16").Value = "  +2.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.761"
$ws.Range("E17").Value = "  -2.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.517.16"
$ws.Range("E18").Value = "  +0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.38"
$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("E20").Value = "  +2.24%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.92"
$ws.Range("E21").Value = "  +6.33%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.39"
$ws.Range("E22").Value = "  +1.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.45"
$ws.Range("E23").Value = "  +1.01%  "

$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("E26").Value = "  +0.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.91"
$ws.Range("E27").Value = "  +4.45%  "

$ws.Range("E28").Value = "  -0.33%  "

$ws.Range("E29").Value = "  +1.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.05"
$ws.Range("E30").Value = "  +1.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.44"
$ws.Range("E31").Value = "  +4.11%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.48"
$ws.Range("E32").Value = "  +0.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.58"
$ws.Range("E33").Value = "  +0.46%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0753"
$ws.Range("E34").Value = "  +2.51%  "

$ws.Range("B35").Value = "ApeXProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.47"
$ws.Range("E35").Value = "  -6.88%  "

$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.45"
$ws.Range("E36").Value = "  +3.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.92"
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("E38").Value = "  +2.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  -1.36%  "

$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("E41").Value = "  -1.07%  "

$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.54"
$ws.Range("E43").Value = "  -3.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.978.30"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.96"
$ws.Range("E46").Value = "  -1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.08"
$ws.Range("E47").Value = "  +6.23%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.723.54"
$ws.Range("E48").Value = "  +0.96%  "

$ws.Range("E49").Value = "  +2.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.84"
$ws.Range("E50").Value = "  -3.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.40"
$ws.Range("E51").Value = "  -1.71%  "
